# Generate Report for Handback
#
# The e8dffba3-... (e2e\e8dffba3-cc15-45e5-9cdc-d6048bd3d3f5.md) file has now
# been handed back in sync with en-US. Update the status / timestamps /
# error-detail columns across the Overview, zh-cn and de-de sheets to
# reflect the successful handback, mirroring the already-"handed back"
# row above it.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the e8dffba3-... file. Columns E (zh-cn) and
# F (de-de) flip from "Ready for handoff" to "Handed back: in sync with
# en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the e8dffba3-... file.
#   C = Status
#   K = Latest Handback DateTime
#   P = Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("K3").Value = "2016-08-12 12:57:05"
$wsZhCn.Range("P3").Value = ""

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the e8dffba3-... file.
#   C = Status
#   K = Latest Handback DateTime
#   P = Error Detail
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("K3").Value = "2016-08-12 12:57:16"
$wsDeDe.Range("P3").Value = ""

# Error Detail column no longer holds long URLs, so its width shrinks
# back down from the wide "fits the long url" size to a narrow default.
$wsZhCn.Columns.Item(16).ColumnWidth = 13
$wsDeDe.Columns.Item(16).ColumnWidth = 13
